# Code refactoring for #8
# Add a new "ratio" column (D) to the results sheet:
#   D2      = B2/C2
#   D3:D5   = B3/C3 filled relatively (Excel stores this as one shared formula)
#   D6      = SUM(D2:D5)/4   (average of the ratios)
# Finally leave the active selection on C11, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("D2").Formula = "=B2/C2"
$ws.Range("D3:D5").Formula = "=B3/C3"
$ws.Range("D6").Formula = "=SUM(D2:D5)/4"

$ws.Range("C11").Select() | Out-Null
